# Update column G ("K") values in the active worksheet per the source data
# regeneration: "use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 1
    3 = 1
    4 = 1
    5 = 0
    6 = 1
    7 = 2
    8 = 1
    9 = 2
    10 = 0
    11 = 1
    12 = 0
    14 = 0
    15 = 3
    16 = 2
    17 = 1
    18 = 0
    19 = 1
    20 = 1
    21 = 2
    22 = 0
    23 = 2
    24 = 2
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    31 = 1
    32 = 0
    33 = 1
    34 = 1
    35 = 0
    36 = 0
    37 = 3
    38 = 0
    39 = 1
    40 = 1
    41 = 0
    42 = 2
    43 = 1
    44 = 2
    45 = 1
    46 = 1
    47 = 0
    48 = 5
    49 = 2
    50 = 1
    51 = 2
    52 = 0
    53 = 1
    54 = 1
    55 = 2
    56 = 1
    57 = 0
    58 = 1
    59 = 3
    60 = 2
    61 = 4
    62 = 1
    63 = 1
    64 = 1
    65 = 0
    66 = 2
    67 = 1
    68 = 0
    69 = 1
    70 = 0
    71 = 1
    72 = 2
    73 = 0
    74 = 1
    75 = 1
    76 = 1
    77 = 0
    78 = 1
    79 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $gValues[$row]
}
